# Fix the subtitle on the title slide: "map filter append" -> "map filter apply"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Replace just the first paragraph's text ("map filter append") in place,
# leaving the second paragraph ("Functions with variable numbers of args")
# untouched.
$chars = $tr.Characters(1, 18)
$chars.Text = "map filter apply"
